$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Entry")

# --- Insert 3 new columns (Month / Day / Year) after column D ---
$ws.Range("E1:G1").EntireColumn.Insert()

# --- New header cells for the inserted columns ---
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"

# --- Row 2 (Transect 1) ---
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 2015
$ws.Range("I2").Value = 42.416289999999996
$ws.Range("J2").Value = -70.907589999999999
$ws.Range("K2").Value = 42.416400000000003
$ws.Range("L2").Value = -70.908019999999993
$ws.Range("P2").Value = 13.89

# --- Row 3 (Transect 4) ---
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 2015
$ws.Range("I3").Value = 42.416200000000003
$ws.Range("J3").Value = -70.906790000000001
$ws.Range("K3").Value = 42.416319999999999
$ws.Range("L3").Value = -70.906360000000006
$ws.Range("P3").Value = 12.78

# --- Row 4 / Row 5 blank-but-styled cells to mirror the template row ---
$ws.Range("E4:G4").HorizontalAlignment = -4108
$ws.Range("E5:G5").HorizontalAlignment = -4108

# --- Data validation range grows by the 3 inserted columns ---
$ws.Range("D74:J1048576").Validation.Delete()
$ws.Range("D74:J1048576").Validation.Add(3, 1, 1, " IN 20, IN 40, OFF 40, OFF 20")

# --- Fix up the view: select J3 on the Data Entry sheet ---
$ws.Activate()
$ws.Range("J3").Select()

# --- Add the "JEKB Changes" log sheet after "Data Entry" ---
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "JEKB Changes"

$ws2.Range("A1").Value = "Date"
$ws2.Range("B1").Value = "Change"
$ws2.Range("E1").Value = "Orig Lat/long"

$ws2.Range("A2").Value = [DateTime]"2017-07-05"
$ws2.Range("A2").NumberFormat = "m/d/yyyy"
$ws2.Range("B2").Value = "Added MDY"
$ws2.Range("E2:H2").Font.Bold = $true
$ws2.Range("E2").Value = "Start Latitude"
$ws2.Range("F2").Value = "Start Longitude"
$ws2.Range("G2").Value = "End Latitude"
$ws2.Range("H2").Value = "End Longitude"

$ws2.Range("A3").Value = [DateTime]"2017-07-05"
$ws2.Range("A3").NumberFormat = "m/d/yyyy"
$ws2.Range("B3").Value = "Fixed lat/long"
$ws2.Range("E3").Value = "42.41.629 N"
$ws2.Range("F3").Value = "70.90.759 W"
$ws2.Range("G3").Value = "42.41.640 N"
$ws2.Range("H3").Value = "70.90.802 W"

$ws2.Range("E4").Value = "42.41.620 N"
$ws2.Range("F4").Value = "70.90.679 W"
$ws2.Range("G4").Value = "42.41.632 N"
$ws2.Range("H4").Value = "70.90.636 W"

$ws2.Range("E3:H4").Select()

# --- Repair the defined name whose sheet-qualified #REF! lost its prefix
#     when the new sheet was inserted ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data Entry!_FilterDatabase") {
        $n.RefersTo = "='Data Entry'!#REF!"
    }
}

$ws.Activate()
